# Edit: add 3 new example rows ("entre 10 et 20", "moins de quinze", "jusqu'à
# cinquante") plus one more "moins de trente euros" / moyen duplicate row, and
# re-sort the data (column A) grouped by "class" (column B), with the original
# "text"/"class" header row moved back up to row 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply a sort on B2:B65 so Excel records the sort state/condition metadata
# (exact final row order is re-asserted explicitly afterwards).
$so = $ws.Sort
$so.SortFields.Clear()
$key = $ws.Range("B2:B65")
$so.SortFields.Add($key)
$so.SetRange($ws.Range("A2:B65"))
$so.Header = 0
$so.Apply()

# Final data set (header + 68 rows), in the exact target row order:
#   row 1       -> header ("text"/"class")
#   rows 2-33   -> class "cher",  sorted by text
#   rows 34-59  -> class "éco",   sorted by text
#   rows 60-65  -> class "moyen", sorted by text
#   rows 66-69  -> newly added example rows, appended at the end
$data = @(
  @("text", "class"),
  @("aucun problème d'argent", "cher"),
  @("bon", "cher"),
  @("cher ", "cher"),
  @("gastronomique ", "cher"),
  @("haut de gamme ", "cher"),
  @("j’ai les moyens ", "cher"),
  @("je m’en fiche ", "cher"),
  @("je m’en fous ", "cher"),
  @("je suis riche ", "cher"),
  @("le meilleur des meilleurs", "cher"),
  @("luxe ", "cher"),
  @("luxueux ", "cher"),
  @("ma bourse est pleine ", "cher"),
  @("moins de cent euros ", "cher"),
  @("plus de vingt ", "cher"),
  @("plutôt bien", "cher"),
  @("plutôt cher", "cher"),
  @("plutôt cher", "cher"),
  @("plutôt pas mal", "cher"),
  @("qualité", "cher"),
  @("standing ", "cher"),
  @("super standing ", "cher"),
  @("très haut de gamme", "cher"),
  @("très haut de gamme", "cher"),
  @("un bon repas ", "cher"),
  @("un plat à vingt ", "cher"),
  @("un restau de qualité", "cher"),
  @("un resto pour le gratin ", "cher"),
  @("un très bon ", "cher"),
  @("un truc chicos ", "cher"),
  @("vraiment bon", "cher"),
  @("vraiment luxueux", "cher"),
  @("à prix modique ", "éco"),
  @("avantageux ", "éco"),
  @("bas de gamme ", "éco"),
  @("bistrot pas cher ", "éco"),
  @("bon marché ", "éco"),
  @("éco", "éco"),
  @("économique ", "éco"),
  @("menu pas cher le midi ", "éco"),
  @("modique ", "éco"),
  @("moins de dix euros ", "éco"),
  @("moins de vingt ", "éco"),
  @("moins de vingt euros ", "éco"),
  @("pas cher", "éco"),
  @("pas cher ", "éco"),
  @("pas cher du tout éco", "éco"),
  @("pas trop cher", "éco"),
  @("peu onéreux", "éco"),
  @("plutôt pas cher", "éco"),
  @("plutôt un troquet ", "éco"),
  @("très bas de gamme", "éco"),
  @("très bon marché ", "éco"),
  @("très économique", "éco"),
  @("très peu cher ", "éco"),
  @("un menu moins de quinze ", "éco"),
  @("un self ", "éco"),
  @("vraiment pas cher", "éco"),
  @("bon rapport qualité prix ", "moyen"),
  @("moins de trente euros ", "moyen"),
  @("moyen de gamme", "moyen"),
  @("ni cher ni pas cher ", "moyen"),
  @("pas un trois étoiles ", "moyen"),
  @("peu importe le prix ", "moyen"),
  @("entre 10 et 20", "éco"),
  @("moins de quinze", "éco"),
  @("jusqu'à cinquante", "cher"),
  @("moins de trente euros ", "moyen")
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 1
  $ws.Cells.Item($r, 1).Value = $data[$i][0]
  $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Match the selection recorded in the saved view (whole header row selected).
$ws.Range("A1:XFD1").Select()
